$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.606.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.509.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.44'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '196.23'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.56%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.14'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.064.23'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '597.09'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.803.84'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.98'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.63'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.518.02'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.983'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.00'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.27'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.87'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.11'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.99'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.28'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.23%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.95%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.08'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.15'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.750.95'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0809'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.86%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.63'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.08'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '492.56'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.02%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.07%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.04%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.45%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.01%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.34'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.44%  '
